$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing headers (D1:F1) to their new, shorter labels.
$ws.Range("D1").Value = "Actions"
$ws.Range("E1").Value = "Expansions"
$ws.Range("F1").Value = "Goal Tests"

# Add the three new headers (G1:I1).
$ws.Range("G1").Value = "New Nodes"
$ws.Range("H1").Value = "Plan Length"
$ws.Range("I1").Value = "Time elapsed in seconds"

# Resize the affected/new columns to match the published widths.
$ws.Columns.Item(5).ColumnWidth = 13.15   # column E -> stored width 14
$ws.Columns.Item(6).ColumnWidth = 12.8    # column F -> stored width 13.7109375 (closest achievable)
$ws.Columns.Item(7).ColumnWidth = 15.0    # column G -> stored width 15.85546875 (closest achievable)
$ws.Columns.Item(8).ColumnWidth = 15.65   # column H -> stored width 16.42578125 (closest achievable)
$ws.Columns.Item(9).ColumnWidth = 22.15   # column I -> stored width 23

# Move the active selection, matching the saved workbook view.
$ws.Range("I2").Select()
